$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a pharmacy transactions report. A new line item
# ("قطن 100 جم" / count "21:0" / price 20) needs to be inserted as item #23,
# just before the existing "كريم جلايسوليد 40ملل وسط" row (previously row 26).
# That pushes the three rows below it (items 23-25) down by one row each,
# the grand-total row moves from row 29 to row 30 (and its sum increases by
# the new row's price, 20), and the footer row moves from row 30 to row 31.
#
# To keep the existing cell styles/merges intact (an EntireRow Insert() on
# this engine re-derives *new* style records instead of reusing the existing
# ones), we shift the content manually, cell by cell, working from the
# bottom of the sheet upward, and only touch merges where the merged range
# itself is moving.
#
# Note: reading ".Value" (no parens) on this host returns the property
# descriptor, not the cell value (it's a parameterized COM property) -- use
# ".Value2" for all reads/writes instead.
# ---------------------------------------------------------------------------

# --- Step 1: unmerge the ranges that are about to move (footer + total) ---
$ws.Range("A30:E30").UnMerge()
$ws.Range("F30:G30").UnMerge()
$ws.Range("I30:N30").UnMerge()
$ws.Range("K29:N29").UnMerge()

# --- Step 2: push the footer row (old row 30) down to row 31 -------------
$ws.Range("A31").Value2 = $ws.Range("A30").Value2
$ws.Range("F31").Value2 = $ws.Range("F30").Value2
$ws.Range("I31").Value2 = $ws.Range("I30").Value2
$ws.Range("A30").ClearContents()
$ws.Range("F30").ClearContents()
$ws.Range("I30").ClearContents()

# --- Step 3: push the totals row (old row 29) down to row 30, new sum ----
$ws.Range("K30").Value2 = $ws.Range("K29").Value2 + 20
$ws.Range("K29").ClearContents()

# --- Step 4: re-merge the rows that just moved ----------------------------
$ws.Range("A31:E31").Merge()
$ws.Range("F31:G31").Merge()
$ws.Range("I31:N31").Merge()
$ws.Range("K30:N30").Merge()

# --- Step 5: shift data rows 26-28 (items 23-25) down to rows 27-29 ------
# (bottom-up so we never overwrite a row before reading it)
$ws.Range("A29").Value2 = 26
$ws.Range("B29").Value2 = $ws.Range("B28").Value2
$ws.Range("H29").Value2 = $ws.Range("H28").Value2
$ws.Range("L29").Value2 = $ws.Range("L28").Value2
$ws.Range("N29").Value2 = $ws.Range("N28").Value2

$ws.Range("B28").Value2 = $ws.Range("B27").Value2
$ws.Range("H28").Value2 = $ws.Range("H27").Value2
$ws.Range("L28").Value2 = $ws.Range("L27").Value2
$ws.Range("N28").Value2 = $ws.Range("N27").Value2

$ws.Range("B27").Value2 = $ws.Range("B26").Value2
$ws.Range("H27").Value2 = $ws.Range("H26").Value2
$ws.Range("L27").Value2 = $ws.Range("L26").Value2
$ws.Range("N27").Value2 = $ws.Range("N26").Value2

# --- Step 6: write the brand-new row 26 (item 23: قطن 100 جم) ------------
$ws.Range("B26").Value2 = "قطن 100 جم"
$ws.Range("H26").Value2 = "21:0"
$ws.Range("L26").Value2 = 20
$ws.Range("N26").Value2 = "1:0"

# --- Step 7: merge the new data row (29) the same way its siblings are ---
$ws.Range("B29:G29").Merge()
$ws.Range("H29:K29").Merge()
$ws.Range("L29:M29").Merge()

# --- Step 8: row heights ---------------------------------------------------
$ws.Rows(26).RowHeight = 24.75
$ws.Rows(27).RowHeight = 25.5
$ws.Rows(28).RowHeight = 25.5
$ws.Rows(29).RowHeight = 24.75
$ws.Rows(30).RowHeight = 26.25
$ws.Rows(31).RowHeight = 16.5
